$wb = $excel.ActiveWorkbook

# Rename "Zipfiles" -> "ZipFiles"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ZipFiles"

# Make the ZipFiles sheet the active tab and set its selection to D19
$ws1.Activate()
$ws1.Range("D19").Select()
